$d = $word.ActiveDocument

# Update the date line (first paragraph) via Find & Replace
$d.Content.Find.Execute("2023-07-26 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-07-27 Thursday", 2) | Out-Null

# Update each table cell by (row, col) position -- avoids ambiguity from duplicate
# source values (e.g. "85x34=2890" appears twice in the table but maps to two
# different targets, so a single global Find/Replace would be unsafe here)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "34×82=2788"
$t.Cell(1,2).Range.Text = "72×17=1224"
$t.Cell(1,3).Range.Text = "35×23=805"
$t.Cell(1,4).Range.Text = "98×89=8722"
$t.Cell(1,5).Range.Text = "58×20=1160"
$t.Cell(2,1).Range.Text = "41×46=1886"
$t.Cell(2,2).Range.Text = "43×99=4257"
$t.Cell(2,3).Range.Text = "40×16=640"
$t.Cell(2,4).Range.Text = "85×16=1360"
$t.Cell(2,5).Range.Text = "89×94=8366"
$t.Cell(3,1).Range.Text = "49×38=1862"
$t.Cell(3,2).Range.Text = "64×64=4096"
$t.Cell(3,3).Range.Text = "76×95=7220"
$t.Cell(3,4).Range.Text = "45×21=945"
$t.Cell(3,5).Range.Text = "73×42=3066"
$t.Cell(4,1).Range.Text = "47×27=1269"
$t.Cell(4,2).Range.Text = "67×58=3886"
$t.Cell(4,3).Range.Text = "12×41=492"
$t.Cell(4,4).Range.Text = "81×50=4050"
$t.Cell(4,5).Range.Text = "65×35=2275"
$t.Cell(5,1).Range.Text = "14×32=448"
$t.Cell(5,2).Range.Text = "32×46=1472"
$t.Cell(5,3).Range.Text = "33×48=1584"
$t.Cell(5,4).Range.Text = "64×37=2368"
$t.Cell(5,5).Range.Text = "24×40=960"
$t.Cell(6,1).Range.Text = "75×22=1650"
$t.Cell(6,2).Range.Text = "13×38=494"
$t.Cell(6,3).Range.Text = "37×29=1073"
$t.Cell(6,4).Range.Text = "17×68=1156"
$t.Cell(6,5).Range.Text = "44×22=968"
$t.Cell(7,1).Range.Text = "43×27=1161"
$t.Cell(7,2).Range.Text = "100×98=9800"
$t.Cell(7,3).Range.Text = "46×99=4554"
$t.Cell(7,4).Range.Text = "93×70=6510"
$t.Cell(7,5).Range.Text = "93×87=8091"
$t.Cell(8,1).Range.Text = "82×95=7790"
$t.Cell(8,2).Range.Text = "48×43=2064"
$t.Cell(8,3).Range.Text = "36×100=3600"
$t.Cell(8,4).Range.Text = "54×73=3942"
$t.Cell(8,5).Range.Text = "78×53=4134"
$t.Cell(9,1).Range.Text = "68×29=1972"
$t.Cell(9,2).Range.Text = "50×90=4500"
$t.Cell(9,3).Range.Text = "58×48=2784"
$t.Cell(9,4).Range.Text = "59×68=4012"
$t.Cell(9,5).Range.Text = "13×30=390"
$t.Cell(10,1).Range.Text = "25×79=1975"
$t.Cell(10,2).Range.Text = "90×56=5040"
$t.Cell(10,3).Range.Text = "42×77=3234"
$t.Cell(10,4).Range.Text = "94×26=2444"
$t.Cell(10,5).Range.Text = "52×18=936"
$t.Cell(11,1).Range.Text = "18×68=1224"
$t.Cell(11,2).Range.Text = "11×37=407"
$t.Cell(11,3).Range.Text = "36×93=3348"
$t.Cell(11,4).Range.Text = "56×13=728"
$t.Cell(11,5).Range.Text = "96×63=6048"
$t.Cell(12,1).Range.Text = "15×63=945"
$t.Cell(12,2).Range.Text = "18×38=684"
$t.Cell(12,3).Range.Text = "46×85=3910"
$t.Cell(12,4).Range.Text = "66×92=6072"
$t.Cell(12,5).Range.Text = "78×78=6084"
$t.Cell(13,1).Range.Text = "66×60=3960"
$t.Cell(13,2).Range.Text = "75×63=4725"
$t.Cell(13,3).Range.Text = "61×86=5246"
$t.Cell(13,4).Range.Text = "19×87=1653"
$t.Cell(13,5).Range.Text = "38×20=760"
$t.Cell(14,1).Range.Text = "35×74=2590"
$t.Cell(14,2).Range.Text = "10×53=530"
$t.Cell(14,3).Range.Text = "72×60=4320"
$t.Cell(14,4).Range.Text = "84×69=5796"
$t.Cell(14,5).Range.Text = "95×31=2945"
$t.Cell(15,1).Range.Text = "85×45=3825"
$t.Cell(15,2).Range.Text = "61×64=3904"
$t.Cell(15,3).Range.Text = "83×61=5063"
$t.Cell(15,4).Range.Text = "41×87=3567"
$t.Cell(15,5).Range.Text = "87×25=2175"
$t.Cell(16,1).Range.Text = "11×32=352"
$t.Cell(16,2).Range.Text = "75×12=900"
$t.Cell(16,3).Range.Text = "60×88=5280"
$t.Cell(16,4).Range.Text = "32×95=3040"
$t.Cell(16,5).Range.Text = "16×40=640"
$t.Cell(17,1).Range.Text = "97×54=5238"
$t.Cell(17,2).Range.Text = "82×82=6724"
$t.Cell(17,3).Range.Text = "67×69=4623"
$t.Cell(17,4).Range.Text = "47×19=893"
$t.Cell(17,5).Range.Text = "73×89=6497"
$t.Cell(18,1).Range.Text = "56×69=3864"
$t.Cell(18,2).Range.Text = "41×43=1763"
$t.Cell(18,3).Range.Text = "42×54=2268"
$t.Cell(18,4).Range.Text = "49×94=4606"
$t.Cell(18,5).Range.Text = "48×96=4608"
$t.Cell(19,1).Range.Text = "69×17=1173"
$t.Cell(19,2).Range.Text = "69×15=1035"
$t.Cell(19,3).Range.Text = "42×75=3150"
$t.Cell(19,4).Range.Text = "63×100=6300"
$t.Cell(19,5).Range.Text = "98×29=2842"
$t.Cell(20,1).Range.Text = "34×16=544"
$t.Cell(20,2).Range.Text = "69×19=1311"
$t.Cell(20,3).Range.Text = "89×76=6764"
$t.Cell(20,4).Range.Text = "80×91=7280"
$t.Cell(20,5).Range.Text = "66×33=2178"
